$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 75956
$ws.Range("E2").Value = 5256
$ws.Range("F2").Value = 5256
$ws.Range("G2").Value = 5858
$ws.Range("H2").Value = 4392
$ws.Range("I2").Value = 4344
$ws.Range("J2").Value = 48
$ws.Range("K2").Value = 59373
$ws.Range("L2").Value = 30884
$ws.Range("M2").Value = 28489
$ws.Range("N2").Value = 28489
$ws.Range("P2").Value = 1360
$ws.Range("Q2").Value = 4486
$ws.Range("R2").Value = -5178
$ws.Range("S2").Value = 1914
$ws.Range("T2").Value = 2908
$ws.Range("U2").Value = 1578
$ws.Range("V2").Value = 12261
$ws.Range("W2").Value = 6.92
$ws.Range("X2").Value = 5.78
$ws.Range("Y2").Value = 17.03
$ws.Range("Z2").Value = 8.130000000000001
$ws.Range("AA2").Value = 108.41
$ws.Range("AB2").Value = 2046.19
$ws.Range("AC2").Value = 16733
$ws.Range("AD2").Value = 10.52
$ws.Range("AE2").Value = 107255
$ws.Range("AF2").Value = 1.64
$ws.Range("AG2").Value = 800
$ws.Range("AH2").Value = 0.45
$ws.Range("AI2").Value = 4.89
$ws.Range("AJ2").Value = 27195083

# Row 3
$ws.Range("D3").Value = 78842
$ws.Range("E3").Value = 5009
$ws.Range("F3").Value = 5009
$ws.Range("G3").Value = 4433
$ws.Range("H3").Value = 3269
$ws.Range("I3").Value = 3269
$ws.Range("K3").Value = 66831
$ws.Range("L3").Value = 35199
$ws.Range("M3").Value = 31632
$ws.Range("N3").Value = 31632
$ws.Range("P3").Value = 1360
$ws.Range("Q3").Value = 6265
$ws.Range("R3").Value = -8750
$ws.Range("S3").Value = 3453
$ws.Range("T3").Value = 5760
$ws.Range("U3").Value = 505
$ws.Range("V3").Value = 15766
$ws.Range("W3").Value = 6.35
$ws.Range("X3").Value = 4.15
$ws.Range("Y3").Value = 10.87
$ws.Range("Z3").Value = 5.18
$ws.Range("AA3").Value = 111.28
$ws.Range("AB3").Value = 2271.59
$ws.Range("AC3").Value = 12020
$ws.Range("AD3").Value = 9.32
$ws.Range("AE3").Value = 119088
$ws.Range("AF3").Value = 0.9399999999999999
$ws.Range("AG3").Value = 1100
$ws.Range("AH3").Value = 0.98
$ws.Range("AI3").Value = 8.94
$ws.Range("AJ3").Value = 27195083

# Row 4
$ws.Range("D4").Value = 75894
$ws.Range("E4").Value = 2627
$ws.Range("F4").Value = 2627
$ws.Range("G4").Value = 2328
$ws.Range("H4").Value = 1307
$ws.Range("I4").Value = 1307
$ws.Range("K4").Value = 70221
$ws.Range("L4").Value = 37635
$ws.Range("M4").Value = 32586
$ws.Range("N4").Value = 32586
$ws.Range("P4").Value = 1360
$ws.Range("Q4").Value = 3546
$ws.Range("R4").Value = -4700
$ws.Range("S4").Value = 3314
$ws.Range("T4").Value = 4383
$ws.Range("U4").Value = -837
$ws.Range("V4").Value = 19616
$ws.Range("W4").Value = 3.46
$ws.Range("X4").Value = 1.72
$ws.Range("Y4").Value = 4.07
$ws.Range("Z4").Value = 1.91
$ws.Range("AA4").Value = 115.49
$ws.Range("AB4").Value = 2346.07
$ws.Range("AC4").Value = 4807
$ws.Range("AD4").Value = 15.19
$ws.Range("AE4").Value = 122680
$ws.Range("AF4").Value = 0.6
$ws.Range("AG4").Value = 1100
$ws.Range("AH4").Value = 1.51
$ws.Range("AI4").Value = 22.35
$ws.Range("AJ4").Value = 27195083

# Row 5
$ws.Range("D5").Value = 74874
$ws.Range("E5").Value = 167
$ws.Range("F5").Value = 167
$ws.Range("G5").Value = -971
$ws.Range("H5").Value = -630
$ws.Range("I5").Value = -630
$ws.Range("K5").Value = 71937
$ws.Range("L5").Value = 40536
$ws.Range("M5").Value = 31401
$ws.Range("N5").Value = 31401
$ws.Range("P5").Value = 1360
$ws.Range("Q5").Value = 727
$ws.Range("R5").Value = -2344
$ws.Range("S5").Value = 3112
$ws.Range("T5").Value = 3081
$ws.Range("U5").Value = -2354
$ws.Range("V5").Value = 22528
$ws.Range("W5").Value = 0.22
$ws.Range("X5").Value = -0.84
$ws.Range("Y5").Value = -1.97
$ws.Range("Z5").Value = -0.89
$ws.Range("AA5").Value = 129.09
$ws.Range("AB5").Value = 2285.12
$ws.Range("AC5").Value = -2317
$ws.Range("AD5").Value = -28.06
$ws.Range("AE5").Value = 118219
$ws.Range("AF5").Value = 0.55
$ws.Range("AG5").Value = 600
$ws.Range("AH5").Value = 0.92
$ws.Range("AI5").Value = -25.3
$ws.Range("AJ5").Value = 27195083

# Row 6
$ws.Range("D6").Value = 78805
$ws.Range("E6").Value = 50
$ws.Range("F6").Value = 50
$ws.Range("G6").Value = -706
$ws.Range("H6").Value = -556
$ws.Range("I6").Value = -556
$ws.Range("K6").Value = 71074
$ws.Range("L6").Value = 40705
$ws.Range("M6").Value = 30369
$ws.Range("N6").Value = 30369
$ws.Range("P6").Value = 1360
$ws.Range("Q6").Value = -1133
$ws.Range("R6").Value = -4639
$ws.Range("S6").Value = -209
$ws.Range("T6").Value = 1740
$ws.Range("U6").Value = -2874
$ws.Range("V6").Value = 22795
$ws.Range("W6").Value = 0.06
$ws.Range("X6").Value = -0.71
$ws.Range("Y6").Value = -1.8
$ws.Range("Z6").Value = -0.78
$ws.Range("AA6").Value = 134.03
$ws.Range("AB6").Value = 2221.64
$ws.Range("AC6").Value = -2043
$ws.Range("AD6").Value = -17.74
$ws.Range("AE6").Value = 114335
$ws.Range("AF6").Value = 0.32
$ws.Range("AG6").Value = 600
$ws.Range("AH6").Value = 1.66
$ws.Range("AI6").Value = -28.68
$ws.Range("AJ6").Value = 27195083

# Row 7
$ws.Range("D7").Value = 75622
$ws.Range("E7").Value = 1295
$ws.Range("G7").Value = 923
$ws.Range("H7").Value = 713
$ws.Range("I7").Value = 714
$ws.Range("K7").Value = 70712
$ws.Range("L7").Value = 39768
$ws.Range("M7").Value = 30945
$ws.Range("N7").Value = 30961
$ws.Range("P7").Value = 1354
$ws.Range("Q7").Value = 4644
$ws.Range("R7").Value = -2644
$ws.Range("S7").Value = -747
$ws.Range("T7").Value = 2323
$ws.Range("U7").Value = 2187
$ws.Range("W7").Value = 1.71
$ws.Range("X7").Value = 0.9399999999999999
$ws.Range("Y7").Value = 2.33
$ws.Range("Z7").Value = 1.01
$ws.Range("AA7").Value = 128.51
$ws.Range("AC7").Value = 2626
$ws.Range("AD7").Value = 19.99
$ws.Range("AE7").Value = 116563
$ws.Range("AF7").Value = 0.45
$ws.Range("AG7").Value = 692
$ws.Range("AH7").Value = 1.32
$ws.Range("AI7").Value = 26.34

# Row 8
$ws.Range("D8").Value = 75673
$ws.Range("E8").Value = 1836
$ws.Range("G8").Value = 1573
$ws.Range("H8").Value = 1217
$ws.Range("I8").Value = 1217
$ws.Range("K8").Value = 71246
$ws.Range("L8").Value = 39406
$ws.Range("M8").Value = 31840
$ws.Range("N8").Value = 31859
$ws.Range("P8").Value = 1354
$ws.Range("Q8").Value = 3618
$ws.Range("R8").Value = -2944
$ws.Range("S8").Value = -244
$ws.Range("T8").Value = 2486
$ws.Range("U8").Value = 1334
$ws.Range("W8").Value = 2.43
$ws.Range("X8").Value = 1.61
$ws.Range("Y8").Value = 3.87
$ws.Range("Z8").Value = 1.72
$ws.Range("AA8").Value = 123.76
$ws.Range("AC8").Value = 4475
$ws.Range("AD8").Value = 9.58
$ws.Range("AE8").Value = 119945
$ws.Range("AF8").Value = 0.36
$ws.Range("AG8").Value = 708
$ws.Range("AH8").Value = 1.65
$ws.Range("AI8").Value = 15.82

# Row 9
$ws.Range("D9").Value = 78826
$ws.Range("E9").Value = 2469
$ws.Range("G9").Value = 2264
$ws.Range("H9").Value = 1732
$ws.Range("I9").Value = 1732
$ws.Range("K9").Value = 73737
$ws.Range("L9").Value = 40404
$ws.Range("M9").Value = 33334
$ws.Range("N9").Value = 33332
$ws.Range("P9").Value = 1354
$ws.Range("Q9").Value = 3992
$ws.Range("R9").Value = -3206
$ws.Range("S9").Value = 56
$ws.Range("T9").Value = 2535
$ws.Range("U9").Value = 1670
$ws.Range("W9").Value = 3.13
$ws.Range("X9").Value = 2.2
$ws.Range("Y9").Value = 5.32
$ws.Range("Z9").Value = 2.39
$ws.Range("AA9").Value = 121.21
$ws.Range("AC9").Value = 6370
$ws.Range("AD9").Value = 6.73
$ws.Range("AE9").Value = 125488
$ws.Range("AF9").Value = 0.34
$ws.Range("AG9").Value = 727
$ws.Range("AH9").Value = 1.7
$ws.Range("AI9").Value = 11.41

# Remove cells that no longer exist in target (O2:O5, J3:J5)
$ws.Range("O2:O5").ClearContents()
$ws.Range("J3:J5").ClearContents()
